# Trim trailing (and leading) whitespace from the NCBI ID values in column B
# so that e.g. "CP016552.1 " becomes "CP016552.1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $cell.Value2.ToString().Trim()
}

# Update the selected/active cell recorded in the sheet view.
$ws.Range("E6").Select()
